$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the asset table: renamed/consolidated categories, new naming
# convention for asset keys, unified Status column, and one added row
# (Grab_Item) for the new Grab Item sound effect.

$ws.Range("A1").Value2 = "Description"
$ws.Range("B1").Value2 = "Category"
$ws.Range("C1").Value2 = "Assets Required"
$ws.Range("D1").Value2 = "Status"

$ws.Range("A2").Value2 = "Game_Ambience"
$ws.Range("B2").Value2 = "Ambience"
$ws.Range("C2").Value2 = "Keyboard noises, loud talking/crowd noises"
$ws.Range("D2").Value2 = "Mixed, implemented in game"

$ws.Range("A3").Value2 = "Walking"
$ws.Range("B3").Value2 = "Ambience"
$ws.Range("C3").Value2 = "Various Footseps"
$ws.Range("D3").Value2 = "Mixed, implemented in game"

$ws.Range("A4").Value2 = "White_Noise"
$ws.Range("B4").Value2 = "Ambience"
$ws.Range("C4").Value2 = "dl from freesound"
$ws.Range("D4").Value2 = "Mixed, implemented in game"

$ws.Range("A5").Value2 = "Game_Over"
$ws.Range("B5").Value2 = "BGM"
$ws.Range("C5").Value2 = "Create/Splice Song from samples/midi"
$ws.Range("D5").Value2 = "Mixed, implemented in game"

$ws.Range("A6").Value2 = "Game_Song"
$ws.Range("B6").Value2 = "BGM"
$ws.Range("C6").Value2 = "Create/Splice Song from samples/midi"
$ws.Range("D6").Value2 = "Mixed, implemented in game"

$ws.Range("A7").Value2 = "Menu"
$ws.Range("B7").Value2 = "BGM"
$ws.Range("C7").Value2 = "Create/Splice Song from samples/midi"
$ws.Range("D7").Value2 = "Mixed, implemented in game"

$ws.Range("A8").Value2 = "Angry"
$ws.Range("B8").Value2 = "Emotion"
$ws.Range("C8").Value2 = "Get some yelling, or protest audio and cut it"
$ws.Range("D8").Value2 = "Mixed, implemented in game"

$ws.Range("A9").Value2 = "Bored"
$ws.Range("B9").Value2 = "Emotion"
$ws.Range("C9").Value2 = "sigh, multiple sighs to randomise"
$ws.Range("D9").Value2 = "Mixed, implemented in game"

$ws.Range("A10").Value2 = "Sad"
$ws.Range("B10").Value2 = "Emotion"
$ws.Range("C10").Value2 = "Get dejected noise, or record myself"
$ws.Range("D10").Value2 = "Mixed, implemented in game"

$ws.Range("A11").Value2 = "Very_Sad"
$ws.Range("B11").Value2 = "Emotion"
$ws.Range("C11").Value2 = "Get crying noises "
$ws.Range("D11").Value2 = "Mixed, implemented in game"

$ws.Range("A12").Value2 = "AC"
$ws.Range("B12").Value2 = "Problem"
$ws.Range("C12").Value2 = "White Noise, Air conditioner on, broken noise?"
$ws.Range("D12").Value2 = "Mixed, implemented in game"

$ws.Range("A13").Value2 = "Apple"
$ws.Range("B13").Value2 = "Problem"
$ws.Range("C13").Value2 = "fridge door open and close"
$ws.Range("D13").Value2 = "Mixed, implemented in game"

$ws.Range("A14").Value2 = "Feed_Fish"
$ws.Range("B14").Value2 = "Problem"
$ws.Range("C14").Value2 = "water sploosh, small items shaking"
$ws.Range("D14").Value2 = "Mixed, implemented in game"

$ws.Range("A15").Value2 = "Pour_Coffee"
$ws.Range("B15").Value2 = "Problem"
$ws.Range("C15").Value2 = "water flowing noise"
$ws.Range("D15").Value2 = "Mixed, implemented in game"

$ws.Range("A16").Value2 = "Router"
$ws.Range("B16").Value2 = "Problem"
$ws.Range("C16").Value2 = "Crash noise, fix noise"
$ws.Range("D16").Value2 = "Mixed, implemented in game"

$ws.Range("A17").Value2 = "Water_Plants"
$ws.Range("B17").Value2 = "Problem"
$ws.Range("C17").Value2 = "dampen a water flowing noise"
$ws.Range("D17").Value2 = "Mixed, implemented in game"

$ws.Range("A18").Value2 = "Correct"
$ws.Range("B18").Value2 = "Sound Effect"
$ws.Range("C18").Value2 = "A ding noise, possibly multiple, chime"
$ws.Range("D18").Value2 = "Mixed, implemented in game"

$ws.Range("A19").Value2 = "Grab_Item"
$ws.Range("B19").Value2 = "Sound Effect"
$ws.Range("C19").Value2 = "swoosh noise, something hitting your palm"
$ws.Range("D19").Value2 = "Mixed, implemented in game"

# Re-sort the final two (newly added) rows alphabetically by asset name,
# matching the sort operation performed in Excel before saving.
$sortRange = $ws.Range("A18:D19")
$sortKey = $ws.Range("A18:A19")
$sortRange.Sort($sortKey, 1)

# Restore the selection/active cell as it was when the workbook was saved.
$ws.Range("C11").Select() | Out-Null
